# Refresh the cryptos.xlsx price/volume snapshot (scheduled GitHub Actions update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeRef, $text) {
    $r = $ws.Range($rangeRef)
    $origStyle = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = $origStyle
}

Set-TextValue "D2" "25.882.48"
$ws.Range("E2").Value = "  +0.38%  "
Set-TextValue "D3" "1.638.93"
$ws.Range("E3").Value = "  +0.94%  "
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  +0.07%  "
Set-TextValue "D5" "215.70"
$ws.Range("E5").Value = "  +0.29%  "
Set-TextValue "D6" "0.508"
$ws.Range("E6").Value = "  -0.29%  "
Set-TextValue "D7" "1.00"
$ws.Range("E7").Value = "  +0.08%  "
Set-TextValue "D8" "0.260"
$ws.Range("E8").Value = "  +1.23%  "
Set-TextValue "D9" "0.0645"
$ws.Range("E9").Value = "  +1.27%  "
Set-TextValue "D10" "20.25"
$ws.Range("E10").Value = "  +4.87%  "
Set-TextValue "D11" "0.0783"
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D12" "1.667.20"
$ws.Range("E12").Value = "  +2.76%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D13" "4.27"
$ws.Range("E13").Value = "  +0.36%  "
Set-TextValue "D14" "1.864.76"
$ws.Range("E14").Value = "  +0.99%  "
Set-TextValue "D15" "0.566"
$ws.Range("E15").Value = "  +2.07%  "
Set-TextValue "D16" "0.0₃0768"
$ws.Range("E16").Value = "  +2.24%  "
Set-TextValue "D17" "63.32"
$ws.Range("E17").Value = "  -0.42%  "
Set-TextValue "D18" "25.887.96"
$ws.Range("E18").Value = "  +0.39%  "
Set-TextValue "D19" "1.00"
$ws.Range("E19").Value = "  +0.13%  "
Set-TextValue "D20" "194.41"
$ws.Range("E20").Value = "  +0.28%  "
Set-TextValue "D21" "4.39"
$ws.Range("E21").Value = "  +1.22%  "
Set-TextValue "D22" "9.96"
$ws.Range("E22").Value = "  +2.01%  "
Set-TextValue "D23" "6.23"
$ws.Range("E23").Value = "  +4.21%  "
Set-TextValue "D24" "1.00"
$ws.Range("E24").Value = "  +0.10%  "
Set-TextValue "D25" "1.75"
$ws.Range("E25").Value = "  -4.21%  "
Set-TextValue "D26" "138.50"
$ws.Range("E26").Value = "  -1.58%  "
Set-TextValue "D27" "0.123"
$ws.Range("E27").Value = "  -3.21%  "
Set-TextValue "D28" "6.85"
$ws.Range("E28").Value = "  +1.81%  "
Set-TextValue "D29" "15.57"
$ws.Range("E29").Value = "  +0.73%  "
Set-TextValue "D30" "1.24"
$ws.Range("E30").Value = "  +0.68%  "
Set-TextValue "D31" "0.0495"
$ws.Range("E31").Value = "  +2.07%  "
Set-TextValue "D32" "3.32"
$ws.Range("E32").Value = "  +0.44%  "
Set-TextValue "D33" "3.25"
$ws.Range("E33").Value = "  +2.23%  "
Set-TextValue "D34" "1.58"
$ws.Range("E34").Value = "  +1.61%  "
Set-TextValue "D35" "2.39"
$ws.Range("E35").Value = "  +1.17%  "
Set-TextValue "D36" "0.909"
$ws.Range("E36").Value = "  +1.83%  "
Set-TextValue "D37" "2.58"
$ws.Range("E37").Value = "  +1.90%  "
Set-TextValue "D38" "0.552"
$ws.Range("E38").Value = "  +0.71%  "
Set-TextValue "D39" "1.126.04"
$ws.Range("E39").Value = "  +0.10%  "
Set-TextValue "D40" "0.0157"
$ws.Range("E40").Value = "  +0.86%  "
Set-TextValue "D41" "1.00"
$ws.Range("E41").Value = "  -0.01%  "
Set-TextValue "D42" "5.50"
$ws.Range("E42").Value = "  -1.36%  "
Set-TextValue "D43" "99.51"
$ws.Range("E43").Value = "  +2.53%  "
Set-TextValue "D44" "0.800"
$ws.Range("E44").Value = "  +0.94%  "
Set-TextValue "D45" "0.0₆0112"
$ws.Range("E45").Value = "  +0.33%  "
Set-TextValue "D46" "55.57"
$ws.Range("E46").Value = "  +1.62%  "
Set-TextValue "D47" "0.424"
$ws.Range("E47").Value = "  -4.12%  "
Set-TextValue "D48" "0.0505"
$ws.Range("E48").Value = "  -0.37%  "
Set-TextValue "D49" "7.66"
$ws.Range("E49").Value = "  +0.75%  "
Set-TextValue "D50" "0.999"
$ws.Range("E50").Value = "  -0.80%  "
Set-TextValue "D51" "1.00"
$ws.Range("E51").Value = "  +0.04%  "
